# edit.ps1
# Applies the "Explosions, Tutorial" update to the Tutorial worksheet:
#  - Rewrites the Exploration description text.
#  - Inserts a new "damage types" tutorial entry after the Attributes section.
#  - Appends a large set of new tutorial entries (Progress, Gear, Crafting,
#    Weapons/Quality/Chanelling/Infusions, Armour/Protection/Sources,
#    Accessories, Skills/Character Skills/Weapon Skills).
#  - Updates sheet dimension / selection accordingly (handled automatically
#    by the engine as cells are written, and via an explicit Select()).

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Tutorial")

$rows = @(
    @{ r = 1; vals = @("Section", "Subsection", "Title", "Min x", "Min y", "Max x", "Max y", "Text") },
    @{ r = 2; vals = @(1, 1, "Exploration", 865, 495, 865, 465, "Beyond The Veil requires constant exploration of your environment, both to survive, and to carry out the Necromancer’s commands. The Explore action allows you to access the map, where you can discover and exploit new destinations.") },
    @{ r = 3; vals = @(2, 1, "The Map", 810, 375, 810, 375, "Characters can discover new regions and travel to existing ones by using the map. The map shows all visited and reachable regions, and is always centred on your camp- the gate.") },
    @{ r = 4; vals = @(2, 2, "Travelling", 808, 860, 808, 30, "Travelling to regions requires grit, shown below the region name. Remember that returning to the gate never requires grit, so you do not need to worry about using up your grit to return home.") },
    @{ r = 5; vals = @(3, 1, "Health", 1020, 950, 10, 35, "Your health bar shows your remaining health, as well as any conditions you are affected by. When it reaches zero, the current character dies. If that character is the Wanderer, then it’s game over.") },
    @{ r = 6; vals = @(3, 2, "Armour", 10, 950, 1020, 35, "Your armour bar shows your current armour protection level. Each piece of armour decreases incoming damage by 5%.") },
    @{ r = 7; vals = @(3, 3, "Adrenaline", 895, 885, 895, 20, "Combat requires focus and keen reflexes. Dealing damage to enemies builds up your adrenaline, which can be used to dash quickly or activate unlocked skills.") },
    @{ r = 8; vals = @(3, 4, "Ammunition", 785, 975, 1010, 10, "This bar shows the remaining ammunition in your weapon. When the magazine is empty, you must reload. Weapons in Beyond The Veil are magical and can recharge ammunition, so you do not need to worry about running out.") },
    @{ r = 9; vals = @(3, 5, "Compass", 990, 975, 750, 10, "You can consume compass charges to find hidden resources, enemies, and other locations of interest. The compass will show the direction of nearby points of interest.") },
    @{ r = 10; vals = @(3, 6, "Enemy", 565, 10, 565, 915, "When facing an enemy, their attributes will show at the top of the screen. This shows their health and armour, as well as a meter that measures your current accuracy. Continuous fire will decrease your accuracy, so fire in short bursts for high accuracy.") },
    @{ r = 11; vals = @(3, 7, "Leaving", 0, 0, 0, 0, "If you have fully explored a region, or you need to escape an overpowering encounter, you can leave by travelling to the ring around the region’s edge.") },
    @{ r = 12; vals = @(3, 8, "Damage", 0, 0, 0, 0, "Attributes") },
    @{ r = 13; vals = @(4, 1, "Characters have 4 attributes- 2 Physical, and 2 Mental:", 645, 320, 645, 645, "Physical") },
    @{ r = 14; vals = @(4, 2, "Fettle- Increases your maximum health in combat. Grit- Determines how far you can travel on the map, as well as movement speed in combat.", 645, 320, 960, 645, "Mental") },
    @{ r = 15; vals = @(4, 3, "Focus- Increases your total compass charges, and your adrenaline recharge rate. Will- Can be used to restore other attributes, and reduces the cooldown of your skills.", 960, 320, 645, 645, "Restoring Attributes") },
    @{ r = 16; vals = @(4, 4, "Attributes can be restored in a number of ways. The simplest two methods are by Sleeping, and by Meditating. Sleeping restores your Fettle and Grit, whilst Meditating restores your Focus and Will. Keeping the fire lit at camp will also restore Fettle and Grit whilst your characters remain there.", 840, 840, 840, 115, "Hunger and Thirst") },
    @{ r = 17; vals = @(4, 5, "Your Hunger and Thirst will decrease over time, eventually killing you if allowed to fully deplete. Eat food and drink water to stave them off. Hot weather will cause you to dehydrate, whilst Cold weather increases your Hunger.", 635, 375, 635, 605, "Finally, there are 4 types of damage in Beyond The Veil. Physical damage is the most basic, is dealt by all weapons, and can reduced by armour. Shatter damage causes damage to armour, increasing the amount of physical damage dealt. Fire damage causes small amounts of direct damage over time. Void damage builds up over time, and upon reaching a critical level, causes large amounts of damage directly to health.") },
    @{ r = 18; vals = @(4, 6, "Progress", 0, 0, 0, 0, "Progress is saved automatically at the dawn of every day. If a character starves, dehydrates, or loses all their health in combat, they will die. If the Wanderer dies, it is game over, and you must either reload at the dawn of the last day, or start a new game. The day can be restarted at any time from the pause menu.") },
    @{ r = 19; vals = @(4, 7, "Gear", 745, 707, 745, 160, "The gear section shows your currently equipped Weapon, Armour, and Accessory. To learn more about each gear type, select the relevant element.") },
    @{ r = 20; vals = @(5, 1, "Crafting", 20, 10, 20, 970, "The resources you find in the world can be used to craft a variety of different items. As you find resources, recipes will be unlocked that can make your journey through the world easier.") },
    @{ r = 21; vals = @(5, 2, "Resources", 0, 0, 0, 0, "Some resources are always dropped by enemies, some can be found scattered throughout the world, and some can only be crafted.") },
    @{ r = 22; vals = @(5, 3, "Radiance", 0, 0, 0, 0, "Salt is dropped by enemies, and can be used to craft one of the most useful resources in the game- Radiance. Radiance can be used to claim a region. Over time, claimed regions will generate either food, water, or resources. These resources are automatically added to your inventory.") },
    @{ r = 23; vals = @(6, 1, "Weapons", 620, 300, 620, 640, "Weapons come in a number of different varieties, all of which with unique firing behaviours.") },
    @{ r = 24; vals = @(6, 2, "Power", 745, 370, 745, 600, "Weapon power affects the damage output of a weapon. Lower power reduces the effectiveness of your weapons, whilst higher power increases it. The current weapon power is shown on this bar.") },
    @{ r = 25; vals = @(6, 3, "Quality", 745, 370, 745, 600, "A weapon’s quality affects it’s maximum power, and is represented by the width of the power bar. Weapons of equal power will have similar damage outputs, regardless of their quality. For this reason it is important to maintain their power level.") },
    @{ r = 26; vals = @(6, 4, "Chanelling", 1190, 700, 565, 275, "Using your weapon slowly drains it of power, thereby decreasing the damage you deal in combat. You can offset this either by using a different weapon, or by channelling essence into it to increase it’s Power. Finding this essence, however, may be easier said than done.") },
    @{ r = 27; vals = @(6, 5, "Infusions", 565, 700, 1190, 275, "Weapons can also be infused with the last words of the dying. These infusions grant larger, passive bonuses to weapons, allowing you to customise them further to your playstyle. Infusing can only be performed using Essence, and new Infusions will replace any existing ones.") },
    @{ r = 28; vals = @(7, 1, "Armour", 505, 305, 505, 660, "Armour can be equipped to reduce incoming damage. All characters have two armour slots for Chest and Head armour respectively.") },
    @{ r = 29; vals = @(7, 2, "Protection", 0, 0, 0, 0, "Armour offers varying levels of protection, with higher quality armour offering more protection. Each point of protection offers a 5% reduction in damage, up to a maximum of 50%.") },
    @{ r = 30; vals = @(7, 3, "Damage", 0, 0, 0, 0, "The Shatter effect causes damage to armour. When enough Shatter damage is taken, armour will lose a point of protection. If a piece of armour loses all it’s protection, it will break, and be lost.") },
    @{ r = 31; vals = @(7, 4, "Sources", 0, 0, 0, 0, "Armour cannot be repaired, and so must be replaced, either through crafting, or found as loot. Instead of a health boost, consider armour an expendable source of protection, destined to be damaged and eventually destroyed.") },
    @{ r = 32; vals = @(8, 1, "Accessories", 0, 0, 0, 0, "Accessories are small trinkets that can further modify the power of your weapon, either by increasing the chance of inflicting conditions, or by increasing the stats of your weapons. Only one accessory can be equipped at a time.") },
    @{ r = 33; vals = @(9, 1, "Skills", 525, 900, 525, 55, "Skills can have a large impact on how a battle plays out. In total there are 4 skill slots available during combat.") },
    @{ r = 34; vals = @(9, 2, "Character Skills", 570, 900, 1020, 55, "2 skills are linked to the character’s class. These skills are unlocked over time, and are generally geared towards survival.") },
    @{ r = 35; vals = @(9, 3, "Weapon Skills", 1020, 900, 570, 55, "The remaining 2 skills are linked to the currently equipped weapon. Weapon skills unlocked by dealing damage with the relevant weapon. Weapon skills have an offensive focus, and are complement the play style of each weapon.") },
    @{ r = 36; vals = @(9, 4, "Cost", 585, 900, 1195, 110, "The cost of a skill is shown above its name. Activating any skill will reset the cooldown period of all your skills, regardless of their cost.") }
)

foreach ($row in $rows) {
    $r = $row.r
    $vals = $row.vals
    for ($i = 0; $i -lt $vals.Length; $i++) {
        $ws.Cells.Item($r, $i + 1).Value = $vals[$i]
    }
}

# Rows 2-19 (column H) carry the "description" style used throughout the
# sheet (style index 1 in the original workbook: vertical-centered text).
# Propagate that style to the three newly inserted rows (17-19) by copying
# it from an existing styled cell in the same column.
$styleSource = $ws.Cells.Item(16, 8).Style
$ws.Cells.Item(17, 8).Style = $styleSource
$ws.Cells.Item(18, 8).Style = $styleSource
$ws.Cells.Item(19, 8).Style = $styleSource

# Match the author's final selection state.
$ws.Range("G37").Select()
